{"js": "async (context) => {\n  // Mapping of old text -> new text, derived from the diff.\n  const replacements = [\n    [\"2025-10-12 Sunday\", \"2025-10-13 Monday\"],\n    [\"937\u00f78=\", \"889\u00f75=\"],\n    [\"148\u00f78=\", \"162\u00f79=\"],\n    [\"997\u00f73=\", \"832\u00f78=\"],\n    [\"376\u00f78=\", \"429\u00f77=\"],\n    [\"579\u00f74=\", \"528\u00f73=\"],\n    [\"542\u00f77=\", \"190\u00f74=\"],\n    [\"860\u00f77=\", \"735\u00f75=\"],\n    [\"379\u00f72=\", \"841\u00f73=\"],\n    [\"205\u00f78=\", \"112\u00f79=\"],\n    [\"755\u00f79=\", \"815\u00f74=\"],\n    [\"121\u00f78=\", \"328\u00f79=\"],\n    [\"245\u00f79=\", \"810\u00f76=\"],\n    [\"384\u00f78=\", \"971\u00f73=\"],\n    [\"950\u00f78=\", \"176\u00f77=\"],\n    [\"251\u00f75=\", \"136\u00f75=\"],\n    [\"654\u00f74=\", \"670\u00f72=\"],\n    [\"359\u00f79=\", \"463\u00f72=\"],\n    [\"399\u00f74=\", \"143\u00f72=\"],\n    [\"445\u00f72=\", \"996\u00f75=\"],\n    [\"622\u00f72=\", \"861\u00f79=\"],\n    [\"472\u00f77=\", \"450\u00f76=\"],\n    [\"985\u00f74=\", \"829\u00f79=\"],\n    [\"686\u00f76=\", \"222\u00f73=\"],\n    [\"902\u00f79=\", \"968\u00f77=\"],\n    [\"525\u00f78=\", \"106\u00f77=\"],\n  ];\n\n  const body = context.document.body;\n\n  for (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n      results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"2025-10-12 Sunday\"; new = \"2025-10-13 Monday\"},\n    @{old = \"937\u00f78=\"; new = \"889\u00f75=\"},\n    @{old = \"148\u00f78=\"; new = \"162\u00f79=\"},\n    @{old = \"997\u00f73=\"; new = \"832\u00f78=\"},\n    @{old = \"376\u00f78=\"; new = \"429\u00f77=\"},\n    @{old = \"579\u00f74=\"; new = \"528\u00f73=\"},\n    @{old = \"542\u00f77=\"; new = \"190\u00f74=\"},\n    @{old = \"860\u00f77=\"; new = \"735\u00f75=\"},\n    @{old = \"379\u00f72=\"; new = \"841\u00f73=\"},\n    @{old = \"205\u00f78=\"; new = \"112\u00f79=\"},\n    @{old = \"755\u00f79=\"; new = \"815\u00f74=\"},\n    @{old = \"121\u00f78=\"; new = \"328\u00f79=\"},\n    @{old = \"245\u00f79=\"; new = \"810\u00f76=\"},\n    @{old = \"384\u00f78=\"; new = \"971\u00f73=\"},\n    @{old = \"950\u00f78=\"; new = \"176\u00f77=\"},\n    @{old = \"251\u00f75=\"; new = \"136\u00f75=\"},\n    @{old = \"654\u00f74=\"; new = \"670\u00f72=\"},\n    @{old = \"359\u00f79=\"; new = \"463\u00f72=\"},\n    @{old = \"399\u00f74=\"; new = \"143\u00f72=\"},\n    @{old = \"445\u00f72=\"; new = \"996\u00f75=\"},\n    @{old = \"622\u00f72=\"; new = \"861\u00f79=\"},\n    @{old = \"472\u00f77=\"; new = \"450\u00f76=\"},\n    @{old = \"985\u00f74=\"; new = \"829\u00f79=\"},\n    @{old = \"686\u00f76=\"; new = \"222\u00f73=\"},\n    @{old = \"902\u00f79=\"; new = \"968\u00f77=\"},\n    @{old = \"525\u00f78=\"; new = \"106\u00f77=\"}\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2\n    )\n}\n"}
